$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append three new questionnaire responses (rows 19-21) ---
# Copy the formatting of the last existing data row (18) down into the
# three new rows so the date column keeps its date/time number format,
# the answer columns keep their borders/alignment, and the trailing
# (empty) columns Q:V keep their style too.
$ws.Range("A18:V18").Copy()
$ws.Range("A19:V21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$data = @(
    @(44911.447951388887, 3, 4, 3, 3, 4, 2, 3, 2, 3, 4, 2, 3, 3, 4, 4),
    @(44911.484768518516, 5, 5, 5, 1, 5, 1, 5, 1, 1, 5, 1, 4, 1, 5, 5),
    @(44914.534548611111, 5, 5, 5, 1, 5, 1, 5, 1, 1, 5, 1, 5, 1, 5, 5)
)

$rowIndex = 19
foreach ($rowValues in $data) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }
    $ws.Rows.Item($rowIndex).RowHeight = 15.75
    $rowIndex = $rowIndex + 1
}

# --- Update the view: scroll so row 16 is the first row under the
#     frozen header, and select the newly entered rows ---
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
[void]$ws.Range("A19:V21").Select()

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9       # xlPaperA4
$ws.PageSetup.Orientation = 1     # xlPortrait
